# CELGNounTradeLongHold.xlsx - record the trade's exit leg (row 2) and
# append the running-principle row (row 3) from the overnight trade run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: this trade was not profitable and is no longer being held.
$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 104.82
$ws.Range("F2").Value = -2.0465377067563892
$ws.Range("G2").Value = $false

# Row 3: updated principle after the trade closed.
$ws.Range("C3").Value = 9795.35

# Re-fit the data columns now that new values were written.
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 7.3333333333333333
$ws.Columns.Item(3).ColumnWidth = 7
$ws.Columns.Item(4).ColumnWidth = 6.5
$ws.Columns.Item(5).ColumnWidth = 6.1666666666666667
$ws.Columns.Item(6).ColumnWidth = 11.6666666666666667
